# Applies the "finished reading excel file and creating onboarding file" edit:
#  - DEVICE_CONFIG: sample value new22 -> new25
#  - BASIC_NETWORK: scroll/selection reset (no data change)
#  - LAYER2: selection reset (no data change, tab no longer active)
#  - DOCKER_IP: trim trailing space off the "DOCKER IP" label, drop the
#               172.16.0.0/16 sample row
#  - NTP: relabel the sample header to DOCKER IP (matching the header style
#         used elsewhere) and add a blank data row
#  - PROXY: add username/password sample values + a http/https list
#           validation on the protocol column; PROXY becomes the active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# DEVICE_CONFIG
# ---------------------------------------------------------------------------
$wsDevice = $wb.Worksheets.Item("DEVICE_CONFIG")
$wsDevice.Range("D2").Value = "new25"
$wsDevice.Activate() | Out-Null
$wsDevice.Range("D8").Select() | Out-Null

# ---------------------------------------------------------------------------
# BASIC_NETWORK
# ---------------------------------------------------------------------------
$wsBasic = $wb.Worksheets.Item("BASIC_NETWORK")
$wsBasic.Activate() | Out-Null
$wsBasic.Range("E2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# ---------------------------------------------------------------------------
# LAYER2
# ---------------------------------------------------------------------------
$wsLayer2 = $wb.Worksheets.Item("LAYER2")
$wsLayer2.Activate() | Out-Null
$wsLayer2.Range("D9").Select() | Out-Null

# ---------------------------------------------------------------------------
# DOCKER_IP
# ---------------------------------------------------------------------------
$wsDocker = $wb.Worksheets.Item("DOCKER_IP")
$wsDocker.Range("B2").Value = "DOCKER IP"
$wsDocker.Range("B3").ClearContents() | Out-Null
$wsDocker.Activate() | Out-Null
$wsDocker.Range("D12").Select() | Out-Null

# ---------------------------------------------------------------------------
# NTP
# ---------------------------------------------------------------------------
$wsNtp = $wb.Worksheets.Item("NTP")
$wsNtp.Range("B1").Value = "DOCKER IP"
# Re-use the same header fill used elsewhere (copy format off A1) instead of
# the one-off red highlight that used to live on this cell.
$wsNtp.Range("A1").Copy() | Out-Null
$wsNtp.Range("B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$wsNtp.Range("A2").Value = 1
$wsNtp.Columns.Item(1).ColumnWidth = 17.92
$wsNtp.Columns.Item(2).ColumnWidth = 28.65
$wsNtp.Activate() | Out-Null
$wsNtp.Range("B9").Select() | Out-Null

# ---------------------------------------------------------------------------
# PROXY
# ---------------------------------------------------------------------------
$wsProxy = $wb.Worksheets.Item("PROXY")
$wsProxy.Range("D2").Value = "username"
$wsProxy.Range("E2").Value = "password"
$wsProxy.Range("C2").Validation.Add(3, 1, 1, '"http, https"') | Out-Null

# PROXY ends up the active / selected tab.
$wsProxy.Activate() | Out-Null
$wsProxy.Range("C14").Select() | Out-Null
